# "DELIM bug" temporary fix — adds a regression-test sheet ("adddelim")
# demonstrating array-like delimited values, and touches a couple of
# sheet selections along the way.

$wb = $excel.ActiveWorkbook

# --- mergeA: selection moves off D2 to J36 (tabSelected moves away too,
#     since it ends up on the new "adddelim" sheet instead) --------------
$mergeA = $wb.Worksheets.Item("mergeA")
$mergeA.Activate()
[void]$mergeA.Range("J36").Select()

# --- promotion: selection grows from the single cell G2 to A1:G2 -------
$promotion = $wb.Worksheets.Item("promotion")
$promotion.Activate()
[void]$promotion.Range("A1:G2").Select()

# --- add the new "adddelim" sheet after the last existing sheet --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "adddelim"

# Fill column-by-column (A, then B, then C) so new shared-string entries
# land in the same order as the reference workbook.
$newSheet.Range("A1").Value = "Index"
$newSheet.Range("A2").Value = 1
$newSheet.Range("A3").Value = 2

$newSheet.Range("B1").Value = "Array_1()"
$newSheet.Range("B2").Value = "a,b,c"
$newSheet.Range("B3").Value = "d,e,f"

$newSheet.Range("C1").Value = "Array_2(Int)"
$newSheet.Range("C2").Value = "1,2,3"
$newSheet.Range("C3").Value = "4,5,6"

# Selection on the new sheet + it becomes the active tab/sheet.
[void]$newSheet.Range("A4:XFD6").Select()
